$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("INS_1")

# Rename shared string "ANNUAL" -> "DAYNITE" wherever it's used as TimeSlice (column B, rows 6..10)
# Update existing row 6 formatting: row height + thick bottom border, and clear style on B6
$ws.Rows.Item(6).RowHeight = 14.65
$ws.Cells.Item(6, 2).Value = "DAYNITE"

# Fill down rows 7-10 replicating row 6's pattern but with different Year values (F column)
$years = @(2035, 2040, 2045, 2050)
for ($i = 0; $i -lt $years.Length; $i++) {
    $r = 7 + $i
    $ws.Rows.Item($r).RowHeight = 14.65
    $ws.Cells.Item($r, 2).Value = "DAYNITE"
    $ws.Cells.Item($r, 3).Value = "ELC_FIN_DEM"
    $ws.Cells.Item($r, 4).Value = "ELC_GRID_RES"
    $ws.Cells.Item($r, 5).Value = "FLO_SHAR"
    $ws.Cells.Item($r, 6).Value = $years[$i]
    $ws.Cells.Item($r, 7).Value = "FX"
    $ws.Cells.Item($r, 8).Value = 0.5

    $ws.Cells.Item($r, 3).Style = $ws.Cells.Item(6, 3).Style
    $ws.Cells.Item($r, 4).Style = $ws.Cells.Item(6, 4).Style
    $ws.Cells.Item($r, 5).Style = $ws.Cells.Item(6, 5).Style
    $ws.Cells.Item($r, 6).Style = $ws.Cells.Item(6, 6).Style
    $ws.Cells.Item($r, 7).Style = $ws.Cells.Item(6, 7).Style
    $ws.Cells.Item($r, 8).Style = $ws.Cells.Item(6, 8).Style
}

# Row 11: blank formatted row (copy style from row 6/10 cells, no values)
$ws.Cells.Item(11, 3).Style = $ws.Cells.Item(6, 3).Style
$ws.Cells.Item(11, 4).Style = $ws.Cells.Item(6, 4).Style
$ws.Cells.Item(11, 5).Style = $ws.Cells.Item(6, 5).Style
$ws.Cells.Item(11, 6).Style = $ws.Cells.Item(6, 6).Style
$ws.Cells.Item(11, 7).Style = $ws.Cells.Item(6, 7).Style
$ws.Cells.Item(11, 8).Style = $ws.Cells.Item(6, 8).Style

$ws.Range("C16").Select()
